# Regenerate the "K" column (G) values for the save_data sheet.
# This reflects the commit: "regen save_data to use K instead of Strike#,
# regen std/mean, calc and write s_vals" — here we just need to write the
# newly-computed K values into column G for the relevant rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 0
    3  = 3
    4  = 0
    5  = 0
    6  = 1
    7  = 2
    8  = 1
    9  = 2
    10 = 1
    11 = 0
    12 = 0
    13 = 1
    14 = 0
    15 = 0
    16 = 2
    17 = 0
    18 = 2
    19 = 1
    20 = 1
    21 = 1
    22 = 1
    23 = 1
    24 = 0
    25 = 1
    26 = 1
    27 = 3
    28 = 0
    30 = 1
    31 = 0
    32 = 1
    33 = 2
    34 = 0
    35 = 0
    36 = 0
    37 = 2
    38 = 0
    39 = 1
    40 = 0
    41 = 0
    42 = 4
    43 = 1
    44 = 1
    45 = 1
    46 = 1
    47 = 0
    48 = 0
    49 = 1
    50 = 2
    51 = 1
    53 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
